# add inital sign out capability
#
# The "Key" sheet maps each bunk to the staff member who signs its
# attendance sheet. Add a third bunk/staff pairing so a third bunk can
# have its attendance initially signed out, and leave the UI focused on
# the "Attendance" sheet (rather than "Key") with its selection parked
# further down the form.

$wb = $excel.ActiveWorkbook

$wsAttendance = $wb.Worksheets.Item("Attendance")
$wsKey = $wb.Worksheets.Item("Key")

# New row of data on the "Key" sheet: a third bunk with its own staff signer.
$wsKey.Range("A9").Value = "Bunk 3"
$wsKey.Range("B9").Value = "Staff Member 5"
$wsKey.Range("C9").Value = "Staff Member 5 ID"

# Park the selection on "Key" where it ends up after the edit.
$wsKey.Activate()
$wsKey.Range("C10").Select()

# "Attendance" becomes the active/selected tab, with its own selection.
$wsAttendance.Activate()
$wsAttendance.Range("K14").Select()
